# The bot's spreadsheet template is being reset for a newer app version:
# the task rows (A2:F5) that held real data get cleared back to an empty
# template (formatting/styles are kept, only values/text are removed),
# and the active selection moves from E6 to A6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the sample/task data rows, keeping their number formats/styles.
$ws.Range("A2:F5").ClearContents() | Out-Null

# Move the saved selection/active cell to A6 (below the now-empty table).
$ws.Range("A6").Select() | Out-Null
